$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 1.67
$ws.Range("I2").Value = 6
$ws.Range("O2").Value = 1.3
$ws.Range("V2").Value = 1.2
$ws.Range("AI2").Value = 80
$ws.Range("H5").Value = 1.55
$ws.Range("L5").Value = 1.28
$ws.Range("T5").Value = 1.8
$ws.Range("W5").Value = 1.17
$ws.Range("X5").Value = 20
$ws.Range("Z5").Value = 11
$ws.Range("AA5").Value = 16
$ws.Range("AD5").Value = 9.800000000000001
$ws.Range("AE5").Value = 17
$ws.Range("AG5").Value = 990
$ws.Range("AH5").Value = 990
$ws.Range("AK5").Value = 95
$ws.Range("AO5").Value = 8
$ws.Range("G7").Value = 1.99
$ws.Range("M7").Value = 1.06
$ws.Range("O7").Value = 1.06
$ws.Range("W7").Value = 2
$ws.Range("F8").Value = 1.41
$ws.Range("G8").Value = 1.42
$ws.Range("I9").Value = 3.75
$ws.Range("K9").Value = 3.65
$ws.Range("N9").Value = 3.15
$ws.Range("W9").Value = 1.64
$ws.Range("AI9").Value = 60
$ws.Range("AJ9").Value = 40
$ws.Range("M10").Value = 1.07
$ws.Range("O10").Value = 1.07
$ws.Range("F11").Value = 6.6
$ws.Range("G11").Value = 7.2
$ws.Range("H11").Value = 1.61
$ws.Range("I11").Value = 1.63
$ws.Range("Q11").Value = 2.14
$ws.Range("V11").Value = 2.58
$ws.Range("Z11").Value = 8.4
$ws.Range("AA11").Value = 15
$ws.Range("AB11").Value = 19.5
$ws.Range("F12").Value = 2.84
$ws.Range("G12").Value = 2.86
$ws.Range("H12").Value = 2.86
$ws.Range("I12").Value = 2.88
$ws.Range("L12").Value = 1.51
$ws.Range("T12").Value = 1.98
$ws.Range("U12").Value = 1.96
$ws.Range("V12").Value = 1.53
$ws.Range("W12").Value = 1.53
$ws.Range("X12").Value = 10.5
$ws.Range("Y12").Value = 9.199999999999999
$ws.Range("Z12").Value = 17
$ws.Range("AC12").Value = 7
$ws.Range("AG12").Value = 13
$ws.Range("AM12").Value = 130
$ws.Range("AN12").Value = 38
$ws.Range("AO12").Value = 38
$ws.Range("I13").Value = 2.1
$ws.Range("G14").Value = 3.8
$ws.Range("N14").Value = 2.7
$ws.Range("Q14").Value = 2.46
$ws.Range("T14").Value = 1.98
$ws.Range("F15").Value = 1.96
$ws.Range("N15").Value = 2.22
$ws.Range("P15").Value = 2.04
$ws.Range("S15").Value = 2.64
$ws.Range("T15").Value = 1.47
$ws.Range("U15").Value = 1.81
$ws.Range("AN15").Value = 19
$ws.Range("AO15").Value = 65
$ws.Range("F16").Value = 2.3
$ws.Range("G16").Value = 2.32
$ws.Range("H16").Value = 3.55
$ws.Range("I16").Value = 3.6
$ws.Range("L16").Value = 1.4
$ws.Range("V16").Value = 1.38
$ws.Range("W16").Value = 1.75
$ws.Range("AK16").Value = 24
$ws.Range("AN16").Value = 19
$ws.Range("AO16").Value = 42
$ws.Range("F17").Value = 1.79
$ws.Range("H17").Value = 5.5
$ws.Range("I17").Value = 5.7
$ws.Range("J17").Value = 3.75
$ws.Range("L17").Value = 1.42
$ws.Range("T17").Value = 1.98
$ws.Range("V17").Value = 1.21
$ws.Range("W17").Value = 2.24
$ws.Range("X17").Value = 13
$ws.Range("Z17").Value = 42
$ws.Range("AA17").Value = 140
$ws.Range("AC17").Value = 8.4
$ws.Range("AD17").Value = 21
$ws.Range("AE17").Value = 80
$ws.Range("AF17").Value = 10
$ws.Range("AG17").Value = 9.800000000000001
$ws.Range("AI17").Value = 85
$ws.Range("AJ17").Value = 18.5
$ws.Range("AK17").Value = 19.5
$ws.Range("AL17").Value = 40
$ws.Range("AM17").Value = 130
$ws.Range("AN17").Value = 12.5
$ws.Range("AO17").Value = 100
$ws.Range("F18").Value = 3.75
$ws.Range("J18").Value = 3.25
$ws.Range("L18").Value = 1.44
$ws.Range("M18").Value = 1.07
$ws.Range("N18").Value = 1.06
$ws.Range("O18").Value = 1.36
$ws.Range("P18").Value = 1.76
$ws.Range("Q18").Value = 2.14
$ws.Range("R18").Value = 1.25
$ws.Range("S18").Value = 3.55
$ws.Range("T18").Value = 1.01
$ws.Range("U18").Value = 1.01
$ws.Range("V18").Value = 1.8
$ws.Range("W18").Value = 1.32
$ws.Range("X18").Value = 14.5
$ws.Range("Y18").Value = 10
$ws.Range("Z18").Value = 15.5
$ws.Range("AA18").Value = 34
$ws.Range("AB18").Value = 15
$ws.Range("AC18").Value = 9
$ws.Range("AD18").Value = 12.5
$ws.Range("AE18").Value = 30
$ws.Range("AF18").Value = 30
$ws.Range("AG18").Value = 19
$ws.Range("AH18").Value = 22
$ws.Range("AI18").Value = 50
$ws.Range("AJ18").Value = 90
$ws.Range("AK18").Value = 60
$ws.Range("AL18").Value = 75
$ws.Range("AM18").Value = 1000
$ws.Range("AN18").Value = 1000
$ws.Range("AO18").Value = 1000
$ws.Range("P19").Value = 2.14
$ws.Range("Q19").Value = 1.83
$ws.Range("AG19").Value = 16
$ws.Range("AH19").Value = 16
$ws.Range("AJ19").Value = 75
$ws.Range("K20").Value = 10.5
